# Updates cryptocurrency price (D) and 1h volume change (E) columns
# to match the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.866.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.771.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9971"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.93%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9962"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.78%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.81%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.27"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07460"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.10%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.38%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.12%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.120"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.298"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.793.99"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001059"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.31%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9962"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.939"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.872.94"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.12%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.06"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.23"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.992.96"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.164"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.84"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.167"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.681"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08995"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.507"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.14%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.32%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.065"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.46%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.61%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06055"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9959"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.846"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.393"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.85%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5959"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.13"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.981"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.145"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06877"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.33%  "
